$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header cells: "_old" -> "_FV2410", "_new" -> "_FV2504" ---------
# Columns A..J hold the "FV2410" (previously "_old") headers, K holds "diff",
# and columns L..U hold the "FV2504" (previously "_new") headers.
$fv2410Headers = @("Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410","Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410")
$fv2504Headers = @("Segmentname_FV2504","Segmentgruppe_FV2504","Segment_FV2504","Datenelement_FV2504","Segment ID_FV2504","Code_FV2504","Qualifier_FV2504","Beschreibung_FV2504","Bedingungsausdruck_FV2504","Bedingung_FV2504")

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2410Headers[$i]
    $ws.Cells.Item(1, $i + 12).Value = $fv2504Headers[$i]
}

# --- Turn the used range into an Excel Table ("Table1") -------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U63"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- Freeze the header row (split/freeze at row 2) -------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
